{"js": "// Replace the date line and the 25 three-digit-by-one-digit multiplication\n// expressions with their new values (see commit: \"Update master to output\n// generated at c8c62b6\"). Every old string is unique in the document, so a\n// straightforward search-and-replace per pair is safe and order-independent.\nconst replacements = [\n  [\"2026-02-17 Tuesday\", \"2026-02-18 Wednesday\"],\n  [\"827\u00d76=\", \"695\u00d74=\"],\n  [\"934\u00d79=\", \"935\u00d75=\"],\n  [\"129\u00d73=\", \"745\u00d78=\"],\n  [\"265\u00d74=\", \"674\u00d78=\"],\n  [\"739\u00d78=\", \"437\u00d72=\"],\n  [\"866\u00d78=\", \"602\u00d72=\"],\n  [\"960\u00d79=\", \"297\u00d75=\"],\n  [\"744\u00d76=\", \"369\u00d77=\"],\n  [\"975\u00d78=\", \"862\u00d79=\"],\n  [\"225\u00d73=\", \"710\u00d75=\"],\n  [\"678\u00d79=\", \"341\u00d75=\"],\n  [\"527\u00d79=\", \"911\u00d79=\"],\n  [\"913\u00d76=\", \"562\u00d76=\"],\n  [\"898\u00d76=\", \"770\u00d76=\"],\n  [\"867\u00d75=\", \"651\u00d74=\"],\n  [\"893\u00d74=\", \"271\u00d75=\"],\n  [\"987\u00d75=\", \"455\u00d72=\"],\n  [\"208\u00d78=\", \"608\u00d75=\"],\n  [\"552\u00d74=\", \"163\u00d79=\"],\n  [\"171\u00d77=\", \"308\u00d74=\"],\n  [\"781\u00d72=\", \"452\u00d72=\"],\n  [\"447\u00d74=\", \"799\u00d76=\"],\n  [\"900\u00d78=\", \"838\u00d74=\"],\n  [\"538\u00d72=\", \"186\u00d72=\"],\n  [\"795\u00d72=\", \"299\u00d72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 three-digit-by-one-digit multiplication\n# expressions with their new values (see commit: \"Update master to output\n# generated at c8c62b6\"). Every old string is unique in the document, so a\n# straightforward Find/Replace-all per pair is safe and order-independent.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-02-17 Tuesday\", \"2026-02-18 Wednesday\"),\n    @(\"827\u00d76=\", \"695\u00d74=\"),\n    @(\"934\u00d79=\", \"935\u00d75=\"),\n    @(\"129\u00d73=\", \"745\u00d78=\"),\n    @(\"265\u00d74=\", \"674\u00d78=\"),\n    @(\"739\u00d78=\", \"437\u00d72=\"),\n    @(\"866\u00d78=\", \"602\u00d72=\"),\n    @(\"960\u00d79=\", \"297\u00d75=\"),\n    @(\"744\u00d76=\", \"369\u00d77=\"),\n    @(\"975\u00d78=\", \"862\u00d79=\"),\n    @(\"225\u00d73=\", \"710\u00d75=\"),\n    @(\"678\u00d79=\", \"341\u00d75=\"),\n    @(\"527\u00d79=\", \"911\u00d79=\"),\n    @(\"913\u00d76=\", \"562\u00d76=\"),\n    @(\"898\u00d76=\", \"770\u00d76=\"),\n    @(\"867\u00d75=\", \"651\u00d74=\"),\n    @(\"893\u00d74=\", \"271\u00d75=\"),\n    @(\"987\u00d75=\", \"455\u00d72=\"),\n    @(\"208\u00d78=\", \"608\u00d75=\"),\n    @(\"552\u00d74=\", \"163\u00d79=\"),\n    @(\"171\u00d77=\", \"308\u00d74=\"),\n    @(\"781\u00d72=\", \"452\u00d72=\"),\n    @(\"447\u00d74=\", \"799\u00d76=\"),\n    @(\"900\u00d78=\", \"838\u00d74=\"),\n    @(\"538\u00d72=\", \"186\u00d72=\"),\n    @(\"795\u00d72=\", \"299\u00d72=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # 2 = wdReplaceAll\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
